$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the effort value in B4 from the number 20 to the text "26.5"
# (the original sheet stores similar values, e.g. B5, as text rather than
# numbers, so we build the new value as text too, without introducing a
# new cell style - using TEXT()+PasteSpecial keeps the cell's format byte
# identical to the original "General" style while recording it as text).
$tmp = $ws.Range("Z1")
$tmp.Formula = '=TEXT(26.5,"0.0")'
$tmp.Copy()
$ws.Range("B4").PasteSpecial(-4163)  # xlPasteValues
$tmp.ClearContents()

# Move the active selection to C10 (reflecting the cell selected after the edit)
$ws.Range("C10").Select()
